$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21: Seed train fermentation ratio -> make Uniform distribution Triangular ---
$ws.Range("F21").Value2 = "Triangular"

# --- Row 31: Centrifuge solids recovery -> make Uniform distribution Triangular ---
$ws.Range("F31").Value2 = "Triangular"

# --- Row 38: Hydrogenation TAL-to-HMP conversion -> make Uniform distribution Triangular,
#     with Midpoint formula (=E38) and recomputed Lower bound (=0.9*E38);
#     also rename reaction attribute from TAL_to_HMTHP to TAL_to_HMP ---
$ws.Range("F38").Value2 = "Triangular"
$ws.Range("G38").Formula = "=0.9*E38"
$ws.Range("H38").Formula = "=E38"
$ws.Range("K38").Value2 = "R401.TAL_to_HMP_rxn.X = x"

# --- Rows 40-45: rename "Dehydration" section to "Etherification & hydrolysis",
#     and rename HMTHP to HMP throughout ---
$ws.Range("A40").Value2 = "Etherification & hydrolysis catalyst Amberlyst70:HMP ratio"
$ws.Range("A41").Value2 = "Etherification & hydrolysis reaction time"
$ws.Range("A42").Value2 = "Etherification & hydrolysis temperature"
$ws.Range("A43").Value2 = "Etherification & hydrolysis HMP-to-PSA conversion"
$ws.Range("K43").Value2 = "R402.HMP_to_PSA_rxn.X = x"
$ws.Range("A44").Value2 = "Etherification & hydrolysis pressure"
$ws.Range("A45").Value2 = "Etherification & hydrolysis spent catalyst Amberlyst70 replacement rate"

# --- Row 47: Ring-opening & hydrolysis PSA-to-KS conversion -> make Uniform distribution
#     Triangular, with Midpoint formula (=E47) and recomputed Lower bound (=0.9*E47) ---
$ws.Range("F47").Value2 = "Triangular"
$ws.Range("G47").Formula = "=0.9*E47"
$ws.Range("H47").Formula = "=E47"

# --- Update sheet view / selection state to match author's final cursor position ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 30
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F47:H47").Select()
